$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells: "Hint" (H1) and "Popup" (I1)
$ws.Cells.Item(1, 8).Value = "Hint"
$ws.Cells.Item(1, 9).Value = "Popup"

# Style them like the other header cells (white Calibri text), but with a new
# distinct dark-navy fill (fgColor FF203864 / bgColor FF333333) rather than the
# purple used by the original A1:D1,F1:G1 headers.
$headerRange = $ws.Range("H1:I1")
$headerRange.Font.Color = 16777215           # white   -> BGR 0xFFFFFF
$headerRange.Interior.Color = 6567968        # FF203864 -> BGR packed 0x00643820
$headerRange.Interior.PatternColor = 3355443 # FF333333 -> BGR packed 0x00333333

# Move the active selection to E7 (matches the authored change)
$ws.Range("E7").Select()
